$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "Total Forecast"
$ws.Range("B28").Value = 400
$ws.Range("C28").Value = 3902
$ws.Range("D28").Value = 482789
$ws.Range("E28").Value = 67326
$ws.Range("F28").Value = 98
$ws.Range("G28").Value = 6371

$ws.Range("K27").Select()
